# The sheet holds one weekly price record per row (Feria Lagunitas de Puerto
# Montt - Apio). This edit adds one new weekly record right after row 410.
#
# Net effect vs. the original file:
#   - rows 1..410 keep their row numbers; row 410's own data is overwritten
#     with the new record's Fecha/Precio values
#   - row 411 (a pre-existing duplicate of the old row 410) is left as-is
#   - a fresh row is opened at 412, populated as a copy of row 411 (so the
#     "Segunda" counterpart row that used to sit at 412 is preserved there),
#     and everything that used to be at rows 412..541 shifts down to 413..542
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Open a new blank row at 412; rows 412..541 shift down to become 413..542.
$ws.Rows.Item(412).Insert()

# The new row 412 is populated as a duplicate of row 411.
$ws.Rows.Item(411).Copy()
$ws.Rows.Item(412).PasteSpecial()
$excel.CutCopyMode = $false

# Row 410 becomes the new weekly record: new Fecha + new min/max/avg price
# and $/Kg; Variedad, Calidad and Volumen are unchanged.
$ws.Cells.Item(410, 4).Value = 45215
$ws.Cells.Item(410, 11).Value = 12000
$ws.Cells.Item(410, 12).Value = 12000
$ws.Cells.Item(410, 13).Value = 12000
$ws.Cells.Item(410, 16).Value = 2000
